$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Discuss Proposal" -> split into two runs: "Discuss Proposa" + "l"
#    (same run formatting on both pieces, just broken across two <w:r>)
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($txt -eq "Discuss Proposal`r") {
        $pStart = $p.Range.Start
        $textLen = $txt.Length - 1   # drop the trailing paragraph mark

        # Shrink the run's text down to "Discuss Proposa"
        $rngAll = $d.Range($pStart, $pStart + $textLen)
        $rngAll.Text = "Discuss Proposa"

        # Re-insert the trailing "l" right after it
        $insPoint = $rngAll.End
        $rngIns = $d.Range($insPoint, $insPoint)
        $rngIns.InsertAfter("l")

        # Toggle formatting on/off to force a distinct run boundary so the
        # "l" doesn't get silently re-merged back into the first run
        $rngL = $d.Range($insPoint, $insPoint + 1)
        $rngL.Font.Bold = $true
        $rngL.Font.Bold = $false

        break
    }
}

# ------------------------------------------------------------------
# 2) "Client Questions" -> "Possible names for software"
#    (this is the first/only occurrence left at this point)
# ------------------------------------------------------------------
$d.Content.Find.Execute("Client Questions", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Possible names for software", 2)

# ------------------------------------------------------------------
# 3) Add a brand-new list paragraph after it, containing "Client Questions"
#    (same list/paragraph formatting, carried over automatically by
#    InsertParagraphAfter)
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Possible names for software`r") {
        $p.Range.InsertParagraphAfter()
        break
    }
}

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Client Questions"
